# Update latest output (run 58)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update cost / unit cost for the single schedule row ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 396.92921775
$wsSchedule.Range("F2").Value = 6.562983097718255

# --- Sheet "Detailed": update Price column (B) and Type column (C) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B21").Value = 36.06011
$wsDetailed.Range("B22").Value = 36.06011
$wsDetailed.Range("B23").Value = 36.06032
$wsDetailed.Range("C23").Value = "historical"
$wsDetailed.Range("B24").Value = 35.88
$wsDetailed.Range("C24").Value = "historical"
$wsDetailed.Range("B25").Value = -18.02293
$wsDetailed.Range("C25").Value = "historical"
$wsDetailed.Range("B26").Value = -18.19512
$wsDetailed.Range("C26").Value = "historical"
$wsDetailed.Range("B27").Value = -14
$wsDetailed.Range("B28").Value = -14
$wsDetailed.Range("B29").Value = -8.12561
$wsDetailed.Range("B30").Value = -4.81333
$wsDetailed.Range("B31").Value = -6.8
$wsDetailed.Range("B32").Value = -6.8
$wsDetailed.Range("B33").Value = 0.7
$wsDetailed.Range("B34").Value = 18.88451
$wsDetailed.Range("B35").Value = 3.18383
$wsDetailed.Range("B36").Value = -0.31005
$wsDetailed.Range("B37").Value = -2.94054
$wsDetailed.Range("B38").Value = -2.88187
$wsDetailed.Range("B39").Value = -2.95405
$wsDetailed.Range("B40").Value = 8.24929
$wsDetailed.Range("B41").Value = 30.48176
$wsDetailed.Range("B42").Value = 30.57917
$wsDetailed.Range("B43").Value = 9.88462
$wsDetailed.Range("B44").Value = 9.81666
$wsDetailed.Range("B45").Value = 22.66264
$wsDetailed.Range("B47").Value = 57.06004
$wsDetailed.Range("B48").Value = 57.06002
$wsDetailed.Range("B49").Value = 56.98
